$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "rules"
$ws.Range("D1").Value = "adaptive_filter"

# Row 2
$ws.Range("D2").Value = "wRLS"
$ws.Range("E2").Value = 40.09541371005085
$ws.Range("F2").Value = 0.8756586931072388
$ws.Range("G2").Value = 32.06318268262779

# Row 3
$ws.Range("D3").Value = "wRLS"
$ws.Range("E3").Value = 40.5412712744984
$ws.Range("F3").Value = 0.8853959427343264
$ws.Range("G3").Value = 32.42764386229245

# Row 4
$ws.Range("D4").Value = "wRLS"
$ws.Range("E4").Value = 39.84280621973048
$ws.Range("F4").Value = 0.8701419039192626
$ws.Range("G4").Value = 31.93067481794548

# Row 5
$ws.Range("D5").Value = "wRLS"
$ws.Range("E5").Value = 39.51545001896316
$ws.Range("F5").Value = 0.8629926497672198
$ws.Range("G5").Value = 31.63527033818688

# Row 6
$ws.Range("D6").Value = "wRLS"
$ws.Range("E6").Value = 39.62291784164025
$ws.Range("F6").Value = 0.8653396796254729
$ws.Range("G6").Value = 31.65533444655503

# Row 7
$ws.Range("D7").Value = "wRLS"
$ws.Range("E7").Value = 39.78073641632546
$ws.Range("F7").Value = 0.8687863383345276
$ws.Range("G7").Value = 31.77269751176194

# Row 8
$ws.Range("D8").Value = "wRLS"
$ws.Range("E8").Value = 39.82722182056744
$ws.Range("F8").Value = 0.8698015504139304
$ws.Range("G8").Value = 31.7780178514535

$wb.Save()
